$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 281
$ws.Cells.Item(281, 1).Value = 9
$ws.Cells.Item(281, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(281, 3).Value = 'Metropolitana'
$ws.Cells.Item(281, 4).Value = 44448
$ws.Cells.Item(281, 5).Value = 13
$ws.Cells.Item(281, 6).Value = 100112024
$ws.Cells.Item(281, 7).Value = 'Choclo'
$ws.Cells.Item(281, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(281, 9).Value = 'Primera'
$ws.Cells.Item(281, 10).Value = 16
$ws.Cells.Item(281, 11).Value = 36000
$ws.Cells.Item(281, 12).Value = 37000
$ws.Cells.Item(281, 13).Value = 36500
$ws.Cells.Item(281, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(281, 15).Value = 'Argentina'
$ws.Cells.Item(281, 16).Value = 730
$ws.Cells.Item(281, 17).Value = 50
$ws.Cells.Item(281, 18).Value = 'Hortaliza'

# Row 282
$ws.Cells.Item(282, 1).Value = 9
$ws.Cells.Item(282, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(282, 3).Value = 'Metropolitana'
$ws.Cells.Item(282, 4).Value = 44448
$ws.Cells.Item(282, 5).Value = 13
$ws.Cells.Item(282, 6).Value = 100112024
$ws.Cells.Item(282, 7).Value = 'Choclo'
$ws.Cells.Item(282, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(282, 9).Value = 'Primera'
$ws.Cells.Item(282, 10).Value = 28
$ws.Cells.Item(282, 11).Value = 33000
$ws.Cells.Item(282, 12).Value = 34000
$ws.Cells.Item(282, 13).Value = 33500
$ws.Cells.Item(282, 14).Value = '$/malla 70 unidades'
$ws.Cells.Item(282, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(282, 16).Value = 479
$ws.Cells.Item(282, 17).Value = 70
$ws.Cells.Item(282, 18).Value = 'Hortaliza'

# Row 283
$ws.Cells.Item(283, 1).Value = 9
$ws.Cells.Item(283, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(283, 3).Value = 'Metropolitana'
$ws.Cells.Item(283, 4).Value = 44167
$ws.Cells.Item(283, 5).Value = 13
$ws.Cells.Item(283, 6).Value = 100112024
$ws.Cells.Item(283, 7).Value = 'Choclo'
$ws.Cells.Item(283, 8).Value = 'Choclero'
$ws.Cells.Item(283, 9).Value = 'Primera'
$ws.Cells.Item(283, 10).Value = 130
$ws.Cells.Item(283, 11).Value = 30000
$ws.Cells.Item(283, 12).Value = 30000
$ws.Cells.Item(283, 13).Value = 30000
$ws.Cells.Item(283, 14).Value = '$/malla 50 unidades'
$ws.Cells.Item(283, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(283, 16).Value = 600
$ws.Cells.Item(283, 17).Value = 50
$ws.Cells.Item(283, 18).Value = 'Hortaliza'

# Row 284
$ws.Cells.Item(284, 1).Value = 9
$ws.Cells.Item(284, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(284, 3).Value = 'Metropolitana'
$ws.Cells.Item(284, 4).Value = 44238
$ws.Cells.Item(284, 5).Value = 13
$ws.Cells.Item(284, 6).Value = 100112024
$ws.Cells.Item(284, 7).Value = 'Choclo'
$ws.Cells.Item(284, 8).Value = 'Choclero'
$ws.Cells.Item(284, 9).Value = 'Primera'
$ws.Cells.Item(284, 10).Value = 5200
$ws.Cells.Item(284, 11).Value = 350
$ws.Cells.Item(284, 12).Value = 350
$ws.Cells.Item(284, 13).Value = 350
$ws.Cells.Item(284, 14).Value = '$/unidad'
$ws.Cells.Item(284, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(284, 16).Value = 350
$ws.Cells.Item(284, 17).Value = 1
$ws.Cells.Item(284, 18).Value = 'Hortaliza'

# Row 285
$ws.Cells.Item(285, 1).Value = 9
$ws.Cells.Item(285, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(285, 3).Value = 'Metropolitana'
$ws.Cells.Item(285, 4).Value = 44238
$ws.Cells.Item(285, 5).Value = 13
$ws.Cells.Item(285, 6).Value = 100112024
$ws.Cells.Item(285, 7).Value = 'Choclo'
$ws.Cells.Item(285, 8).Value = 'Choclero'
$ws.Cells.Item(285, 9).Value = 'Primera'
$ws.Cells.Item(285, 10).Value = 7000
$ws.Cells.Item(285, 11).Value = 350
$ws.Cells.Item(285, 12).Value = 400
$ws.Cells.Item(285, 13).Value = 375
$ws.Cells.Item(285, 14).Value = '$/unidad'
$ws.Cells.Item(285, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(285, 16).Value = 375
$ws.Cells.Item(285, 17).Value = 1
$ws.Cells.Item(285, 18).Value = 'Hortaliza'

# Row 286
$ws.Cells.Item(286, 1).Value = 9
$ws.Cells.Item(286, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(286, 3).Value = 'Metropolitana'
$ws.Cells.Item(286, 4).Value = 44238
$ws.Cells.Item(286, 5).Value = 13
$ws.Cells.Item(286, 6).Value = 100112024
$ws.Cells.Item(286, 7).Value = 'Choclo'
$ws.Cells.Item(286, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(286, 9).Value = 'Primera'
$ws.Cells.Item(286, 10).Value = 1600
$ws.Cells.Item(286, 11).Value = 150
$ws.Cells.Item(286, 12).Value = 200
$ws.Cells.Item(286, 13).Value = 175
$ws.Cells.Item(286, 14).Value = '$/unidad'
$ws.Cells.Item(286, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(286, 16).Value = 175
$ws.Cells.Item(286, 17).Value = 1
$ws.Cells.Item(286, 18).Value = 'Hortaliza'
$ws.Cells.Item(286, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 287
$ws.Cells.Item(287, 1).Value = 9
$ws.Cells.Item(287, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(287, 3).Value = 'Metropolitana'
$ws.Cells.Item(287, 4).Value = 44399
$ws.Cells.Item(287, 5).Value = 13
$ws.Cells.Item(287, 6).Value = 100112024
$ws.Cells.Item(287, 7).Value = 'Choclo'
$ws.Cells.Item(287, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(287, 9).Value = 'Primera'
$ws.Cells.Item(287, 10).Value = 34
$ws.Cells.Item(287, 11).Value = 18000
$ws.Cells.Item(287, 12).Value = 19000
$ws.Cells.Item(287, 13).Value = 18500
$ws.Cells.Item(287, 14).Value = '$/malla 60 unidades'
$ws.Cells.Item(287, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(287, 16).Value = 308
$ws.Cells.Item(287, 17).Value = 60
$ws.Cells.Item(287, 18).Value = 'Hortaliza'
$ws.Cells.Item(287, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 288
$ws.Cells.Item(288, 1).Value = 9
$ws.Cells.Item(288, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(288, 3).Value = 'Metropolitana'
$ws.Cells.Item(288, 4).Value = 44399
$ws.Cells.Item(288, 5).Value = 13
$ws.Cells.Item(288, 6).Value = 100112024
$ws.Cells.Item(288, 7).Value = 'Choclo'
$ws.Cells.Item(288, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(288, 9).Value = 'Primera'
$ws.Cells.Item(288, 10).Value = 38
$ws.Cells.Item(288, 11).Value = 19000
$ws.Cells.Item(288, 12).Value = 20000
$ws.Cells.Item(288, 13).Value = 19342
$ws.Cells.Item(288, 14).Value = '$/malla 70 unidades'
$ws.Cells.Item(288, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(288, 16).Value = 276
$ws.Cells.Item(288, 17).Value = 70
$ws.Cells.Item(288, 18).Value = 'Hortaliza'
$ws.Cells.Item(288, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
